$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns: "MaxHP" (before HP) and "MaxMP" (before MP) ---
$ws.Columns("E:E").Insert()
$ws.Columns("G:G").Insert()

# --- Relocate the "Satiety" column so it sits right after MP (column I) ---
# Satiety now lives in column T (20) after the two inserts above. Open up a
# slot at I by inserting a blank column there (Satiety shifts from T to U),
# copy the values across, then drop the now-empty original column.
$ws.Columns("I:I").Insert()
for ($r = 1; $r -le 18; $r++) {
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($r, 21).Value2
}
$ws.Columns("U:U").Delete()

# --- Header text for the two brand-new columns ---
$ws.Range("E1").Value = "MaxHP"
$ws.Range("G1").Value = "MaxMP"

# --- Fill the new MaxHP / MaxMP columns with 0 for every data row ---
$ws.Range("E2:E18").Value = 0
$ws.Range("G2:G18").Value = 0

# --- Fix up "Type" values (column C) for a handful of equipment items ---
$ws.Range("C10").Value = 3
$ws.Range("C11").Value = 3
$ws.Range("C12").Value = 3
$ws.Range("C13").Value = 3
$ws.Range("C14").Value = 3
$ws.Range("C15").Value = 2

# --- Keep the hidden AutoFilter-database name in sync with the new extent ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ItemData!_FilterDatabase") {
        $n.RefersTo = "=ItemData!`$A`$1:`$U`$18"
    }
}

# --- Update the saved selection ---
$ws.Range("D1").Select()
